# Horarios actualizados Linea 141 - 964
# Refresh the scraped schedule data (Hora_Scrap/Hora_Llegada/Linea/Minutos/Parada)
# for all three sheets: LP1912, LP1912-215 and 6203-6173.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- LP1912 ---
$ws1.Cells.Item(2, 1).Value = 'Última actualización: 10:41:48'
$ws1.Cells.Item(3, 1).Value = 'Total filas: 130'
$ws1.Cells.Item(16, 1).Value = '05:44:02'
$ws1.Cells.Item(16, 3).Value = '17X38_ROMERO'
$ws1.Cells.Item(16, 4).Value = 56
$ws1.Cells.Item(17, 1).Value = '06:38:54'
$ws1.Cells.Item(17, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(17, 4).Value = 2
$ws1.Cells.Item(41, 1).Value = '06:38:54'
$ws1.Cells.Item(41, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(41, 4).Value = 82
$ws1.Cells.Item(42, 1).Value = '07:52:32'
$ws1.Cells.Item(42, 3).Value = '17_ROMERO'
$ws1.Cells.Item(42, 4).Value = 8
$ws1.Cells.Item(49, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(50, 3).Value = '15_ABASTO'
$ws1.Cells.Item(51, 3).Value = '15_ABASTO'
$ws1.Cells.Item(52, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(85, 1).Value = '08:40:59'
$ws1.Cells.Item(85, 3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(85, 4).Value = 62
$ws1.Cells.Item(86, 1).Value = '08:30:14'
$ws1.Cells.Item(86, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(86, 4).Value = 72
$ws1.Cells.Item(105, 1).Value = '10:41:48'
$ws1.Cells.Item(105, 4).Value = 3
$ws1.Cells.Item(106, 1).Value = '10:41:48'
$ws1.Cells.Item(106, 4).Value = 5
$ws1.Cells.Item(108, 1).Value = '10:41:48'
$ws1.Cells.Item(108, 2).Value = '10:55'
$ws1.Cells.Item(108, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(108, 4).Value = 14
$ws1.Cells.Item(109, 1).Value = '10:41:48'
$ws1.Cells.Item(109, 2).Value = '10:56'
$ws1.Cells.Item(109, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(109, 4).Value = 15
$ws1.Cells.Item(110, 1).Value = '10:41:48'
$ws1.Cells.Item(110, 2).Value = '10:59'
$ws1.Cells.Item(110, 4).Value = 18
$ws1.Cells.Item(111, 1).Value = '09:23:52'
$ws1.Cells.Item(111, 3).Value = '10_OLMOS'
$ws1.Cells.Item(111, 4).Value = 98
$ws1.Cells.Item(112, 1).Value = '10:41:48'
$ws1.Cells.Item(112, 2).Value = '11:01'
$ws1.Cells.Item(112, 3).Value = '81_EL PELIGRO'
$ws1.Cells.Item(112, 4).Value = 20
$ws1.Cells.Item(113, 2).Value = '11:03'
$ws1.Cells.Item(113, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(113, 4).Value = 56
$ws1.Cells.Item(114, 1).Value = '10:41:48'
$ws1.Cells.Item(114, 2).Value = '11:07'
$ws1.Cells.Item(114, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(114, 4).Value = 26
$ws1.Cells.Item(115, 1).Value = '10:41:48'
$ws1.Cells.Item(115, 2).Value = '11:10'
$ws1.Cells.Item(115, 3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(115, 4).Value = 29
$ws1.Cells.Item(116, 1).Value = '10:41:48'
$ws1.Cells.Item(116, 2).Value = '11:14'
$ws1.Cells.Item(116, 3).Value = '14_ABASTO'
$ws1.Cells.Item(116, 4).Value = 33
$ws1.Cells.Item(117, 1).Value = '10:41:48'
$ws1.Cells.Item(117, 2).Value = '11:15'
$ws1.Cells.Item(117, 3).Value = '15X38_ABASTO'
$ws1.Cells.Item(117, 4).Value = 34
$ws1.Cells.Item(118, 1).Value = '09:23:52'
$ws1.Cells.Item(118, 2).Value = '11:19'
$ws1.Cells.Item(118, 4).Value = 116
$ws1.Cells.Item(119, 2).Value = '11:21'
$ws1.Cells.Item(119, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(119, 4).Value = 74
$ws1.Cells.Item(120, 1).Value = '10:41:48'
$ws1.Cells.Item(120, 2).Value = '11:25'
$ws1.Cells.Item(120, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(120, 4).Value = 44
$ws1.Cells.Item(121, 1).Value = '10:41:48'
$ws1.Cells.Item(121, 2).Value = '11:29'
$ws1.Cells.Item(121, 3).Value = '10_OLMOS'
$ws1.Cells.Item(121, 4).Value = 48
$ws1.Cells.Item(122, 1).Value = '10:41:48'
$ws1.Cells.Item(122, 2).Value = '11:30'
$ws1.Cells.Item(122, 3).Value = '215C_EL PATO'
$ws1.Cells.Item(122, 4).Value = 49
$ws1.Cells.Item(123, 1).Value = '10:41:48'
$ws1.Cells.Item(123, 2).Value = '11:41'
$ws1.Cells.Item(123, 3).Value = '215B_EL PATO'
$ws1.Cells.Item(123, 4).Value = 60
$ws1.Cells.Item(124, 1).Value = '10:41:48'
$ws1.Cells.Item(124, 2).Value = '11:45'
$ws1.Cells.Item(124, 3).Value = '15X38_ABASTO'
$ws1.Cells.Item(124, 4).Value = 64
$ws1.Cells.Item(125, 1).Value = '10:41:48'
$ws1.Cells.Item(125, 2).Value = '11:49'
$ws1.Cells.Item(125, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(125, 4).Value = 68
$ws1.Cells.Item(125, 5).Value = 'LP1912'
$ws1.Cells.Item(126, 1).Value = '10:07:51'
$ws1.Cells.Item(126, 2).Value = '11:52'
$ws1.Cells.Item(126, 3).Value = '225_GOMEZ'
$ws1.Cells.Item(126, 4).Value = 105
$ws1.Cells.Item(126, 5).Value = 'LP1912'
$ws1.Cells.Item(127, 1).Value = '10:41:48'
$ws1.Cells.Item(127, 2).Value = '11:53'
$ws1.Cells.Item(127, 3).Value = '225_GOMEZ'
$ws1.Cells.Item(127, 4).Value = 72
$ws1.Cells.Item(127, 5).Value = 'LP1912'
$ws1.Cells.Item(128, 1).Value = '10:41:48'
$ws1.Cells.Item(128, 2).Value = '11:58'
$ws1.Cells.Item(128, 3).Value = '17_ROMERO'
$ws1.Cells.Item(128, 4).Value = 77
$ws1.Cells.Item(128, 5).Value = 'LP1912'
$ws1.Cells.Item(129, 1).Value = '10:41:48'
$ws1.Cells.Item(129, 2).Value = '12:05'
$ws1.Cells.Item(129, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(129, 4).Value = 84
$ws1.Cells.Item(129, 5).Value = 'LP1912'
$ws1.Cells.Item(130, 1).Value = '10:41:48'
$ws1.Cells.Item(130, 2).Value = '12:10'
$ws1.Cells.Item(130, 3).Value = '15_ABASTO'
$ws1.Cells.Item(130, 4).Value = 89
$ws1.Cells.Item(130, 5).Value = 'LP1912'
$ws1.Cells.Item(131, 1).Value = '10:41:48'
$ws1.Cells.Item(131, 2).Value = '12:10'
$ws1.Cells.Item(131, 3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(131, 4).Value = 89
$ws1.Cells.Item(131, 5).Value = 'LP1912'
$ws1.Cells.Item(132, 1).Value = '10:41:48'
$ws1.Cells.Item(132, 2).Value = '12:21'
$ws1.Cells.Item(132, 3).Value = '215C_EL PATO'
$ws1.Cells.Item(132, 4).Value = 100
$ws1.Cells.Item(132, 5).Value = 'LP1912'
$ws1.Cells.Item(133, 1).Value = '10:41:48'
$ws1.Cells.Item(133, 2).Value = '12:32'
$ws1.Cells.Item(133, 3).Value = '14_ABASTO'
$ws1.Cells.Item(133, 4).Value = 111
$ws1.Cells.Item(133, 5).Value = 'LP1912'
$ws1.Cells.Item(134, 1).Value = '10:41:48'
$ws1.Cells.Item(134, 2).Value = '12:34'
$ws1.Cells.Item(134, 3).Value = '15_ABASTO'
$ws1.Cells.Item(134, 4).Value = 113
$ws1.Cells.Item(134, 5).Value = 'LP1912'
$ws1.Cells.Item(135, 1).Value = '10:41:48'
$ws1.Cells.Item(135, 2).Value = '12:36'
$ws1.Cells.Item(135, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(135, 4).Value = 115
$ws1.Cells.Item(135, 5).Value = 'LP1912'

# --- LP1912-215 ---
$ws2.Cells.Item(2, 1).Value = 'Última actualización: 10:41:48'
$ws2.Cells.Item(3, 1).Value = 'Total filas: 18'
$ws2.Cells.Item(21, 1).Value = '10:41:48'
$ws2.Cells.Item(21, 4).Value = 49
$ws2.Cells.Item(22, 1).Value = '10:41:48'
$ws2.Cells.Item(22, 4).Value = 60
$ws2.Cells.Item(23, 1).Value = '10:41:48'
$ws2.Cells.Item(23, 2).Value = '12:21'
$ws2.Cells.Item(23, 3).Value = '215C_EL PATO'
$ws2.Cells.Item(23, 4).Value = 100
$ws2.Cells.Item(23, 5).Value = 'LP1912'

# --- 6203-6173 ---
$ws3.Cells.Item(2, 1).Value = 'Última actualización: 10:41:48'
$ws3.Cells.Item(22, 1).Value = '10:41:48'
$ws3.Cells.Item(22, 4).Value = 44

Write-Host "Updated 160 cells across LP1912, LP1912-215 and 6203-6173."
